$wb = $excel.ActiveWorkbook

# --- Sheet 1: Transfer Time (s) ---
$ws1 = $wb.Worksheets.Item("Transfer Time (s)")
$ws1.Range("D4").Value = 0.008527708172798156
$ws1.Range("E4").Value = 0.003864228111583123
$ws1.Range("D5").Value = 0.01721251964569092
$ws1.Range("E5").Value = 0.01080990045973373
$ws1.Range("D6").Value = 0.1141218423843384
$ws1.Range("E6").Value = 0.02248933364168916
$ws1.Range("D7").Value = 0.9106026887893677

# --- Sheet 2: Throughput (bps) ---
$ws2 = $wb.Worksheets.Item("Throughput (bps)")
$ws2.Range("D4").Value = 1406182.728407707
$ws2.Range("E4").Value = 515828.0395166186
$ws2.Range("D5").Value = 6991508.194169387
$ws2.Range("E5").Value = 2265222.321775638
$ws2.Range("D6").Value = 9829286.040962379
$ws2.Range("E6").Value = 1864039.675413585
$ws2.Range("D7").Value = 11568075.78619621

# --- Sheet 3: Overhead Ratio ---
$ws3 = $wb.Worksheets.Item("Overhead Ratio")
$ws3.Range("D4").Value = 1.020264990234375
$ws3.Range("E4").Value = 0.00004171876828768403
$ws3.Range("D5").Value = 1.0025732421875
$ws3.Range("D6").Value = 1.000741004943848
$ws3.Range("D7").Value = 1.000568580627441
